# [APP]: refactor movie script object
#
# The template previously referenced the root data object as `data`, with
# paths like `data.script.title` or `data.character_list`. The data object
# was refactored so the root is now called `script` directly, so every
# `{{...}}` placeholder that used to read `data.X` must now read `script.X`
# (and `data.script.X` simply collapses to `script.X`).
#
# This touches every templating placeholder in the body (title/author/
# genre/date/plot summary, the locations/props/character-list/director/
# producer/accomplishments/other-scripts/produced-movies/companies loops)
# plus the "About {{...}}" heading, and the header title gains the word
# "script" ("Movie synopsis" -> "Movie script synopsis").

$d = $word.ActiveDocument

# Collapse "data.script." -> "script." first (handles the common
# `{{data.script.<field>}}` placeholders so they don't become
# `script.script.<field>`), then mop up any remaining standalone
# "data." root references (`data.character_list`, `data.directors`,
# `data.accomplishments`, `data.other_scripts`, `data.produced_movies`,
# `data.companies_worked_with`, ...) to the new `script.` root.
$d.Content.Find.Execute("data.script.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "script.", 2) | Out-Null
$d.Content.Find.Execute("data.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "script.", 2) | Out-Null

# The document title in the first-page header gains the word "script":
# "Movie synopsis" -> "Movie script synopsis".
$section = $d.Sections.Item(1)
$firstPageHeader = $section.Headers.Item(2)
$firstPageHeader.Range.Find.Execute("Movie synopsis", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "Movie script synopsis", 2) | Out-Null
